$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the "Send Confirmation Email" column
$ws.Range("E1").Value = "Send Confirmation Email"

# Update employee email addresses (investor access can suppress emails now)
$ws.Range("C2").Value = "emp10@mycompany.com"
$ws.Range("C3").Value = "emp20@mycompany.com"

# Set the new "Send Confirmation Email" flags per investor
$ws.Range("E2").Value = "No"
$ws.Range("E3").Value = "Yes"

# Remove the hyperlink for the first investor's email, keep & refresh the second
$ws.Range("C2").Hyperlinks.Delete()

$h = $ws.Range("C3").Hyperlinks.Item(1)
$h.Address = "mailto:emp20@mycompany.com"
$h.TextToDisplay = "emp20@mycompany.com"

# Update selection to match the saved state
$ws.Range("E3").Select()
